$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 13:52"

# Updated province/ciudad COVID data, re-sorted descending by Casos totales (col B)
$data = @(
    @("Madrid", 38723, 17322, 16265, 5136),
    @("Cataluña", 26824, 10088, 13976, 2760),
    @("Galicia", 6151, 610, 5300, 241),
    @("Bizkaia/Vizcaya", 4603, 3568, 4694, 266),
    @("Ciudad Real", 4125, 1259, 8211, 364),
    @("Valencia/Valencia", 3788, 561, 2929, 298),
    @("Navarra", 3231, 380, 2664, 187),
    @("Castilla-La Mancha", 2780, 71, 2446, 263),
    @("Albacete", 2751, 1259, 8211, 241),
    @("La Rioja", 2719, 964, 1614, 141),
    @("Araba/Alava", 2675, 3568, 4694, 201),
    @("Alacant/Alicante", 2673, 284, 2120, 269),
    @("Zaragoza", 2520, 408, 1895, 217),
    @("Toledo", 2283, 1259, 8211, 307),
    @("A Coruña", 1969, 333, 1788, 67),
    @("Malaga", 1916, 258, 1537, 121),
    @("Salamanca", 1719, 402, 1108, 209),
    @("Sevilla", 1663, 95, 1471, 97),
    @("Asturias", 1646, 222, 1338, 86),
    @("Pontevedra", 1536, 333, 1411, 30),
    @("Gipuzkoa/Guipuzcoa", 1532, 3568, 4694, 81),
    @("Granada", 1501, 129, 1258, 114),
    @("Cantabria", 1483, 118, 1288, 77),
    @("Valladolid", 1445, 487, 813, 145),
    @("Caceres", 1379, 87, 1097, 195),
    @("Leon", 1315, 498, 631, 186),
    @("Murcia", 1259, 156, 1035, 68),
    @("Segovia", 1252, 361, 782, 109),
    @("Cordoba", 1046, 84, 925, 37),
    @("Burgos", 1003, 360, 530, 113),
    @("Tenerife", 971, 157, 1649, 51),
    @("Jaen", 934, 61, 821, 52),
    @("Aragon", 907, 29, 838, 40),
    @("Guadalajara", 873, 1259, 8211, 124),
    @("Castello/Castellon", 869, 85, 714, 70),
    @("Cadiz", 862, 104, 727, 31),
    @("Soria", 824, 164, 596, 64),
    @("Ourense", 751, 333, 660, 22),
    @("Avila", 732, 225, 428, 79),
    @("Badajoz", 689, 136, 520, 33),
    @("Lugo", 586, 333, 520, 11),
    @("Cuenca", 570, 1259, 8211, 96),
    @("Palencia", 486, 114, 333, 39),
    @("Gran Canaria", 467, 157, 1649, 22),
    @("Huesca", 417, 57, 323, 37),
    @("Teruel", 380, 93, 258, 29),
    @("Almeria", 371, 49, 299, 23),
    @("Zamora", 340, 100, 202, 38),
    @("Huelva", 288, 18, 254, 16),
    @("Mallorca", 210, 18, 194, 12),
    @("Melilla", 86, 11, 73, 2),
    @("Ceuta", 83, 11, 69, 3),
    @("La Palma", 69, 157, 1649, 4),
    @("Lanzarote", 69, 157, 1649, 3),
    @("Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena", 58, 0, 58, 3),
    @("Fuerteventura", 35, 157, 1649, 0),
    @("Ibiza", 21, 18, 20, 1),
    @("Menorca", 15, 18, 13, 0),
    @("La Gomera", 8, 157, 1649, 0),
    @("Arroyo de la Luz", 7, 0, 7, 0),
    @("El Hierro", 3, 157, 1649, 0),
    @("Formentera", 0, 10, 0, 8),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 4 + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}
